$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B30").NumberFormat = "@"

$ws.Range("A30").Value = "Nishant"
$ws.Range("B30").Value = "9880188877"
$ws.Range("C30").Value = "Idk"
$ws.Range("D30").Value = "2025-10-01 18:06:20"
